$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.473.63'
$ws.Range("E2").Value = '  -2.61%  '
$ws.Range("D3").Value = '1.777.86'
$ws.Range("E3").Value = '  -2.89%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '''229.62'
$ws.Range("E5").Value = '  -1.96%  '
$ws.Range("D6").Value = '''0.5858'
$ws.Range("E6").Value = '  -2.32%  '
$ws.Range("D7").Value = '''1.002'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = '''0.2741'
$ws.Range("E8").Value = '  -0.65%  '
$ws.Range("D9").Value = '''23.21'
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("D10").Value = '''0.06686'
$ws.Range("E10").Value = '  -4.32%  '
$ws.Range("D11").Value = '''0.07529'
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("D12").Value = '1.781.48'
$ws.Range("E12").Value = '  -2.75%  '
$ws.Range("D13").Value = '''4.746'
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("D14").Value = '''0.6071'
$ws.Range("E14").Value = '  -3.31%  '
$ws.Range("D15").Value = '2.019.13'
$ws.Range("E15").Value = '  -2.81%  '
$ws.Range("D16").Value = '''74.65'
$ws.Range("E16").Value = '  -4.82%  '
$ws.Range("D17").Value = '''0.000008632'
$ws.Range("E17").Value = '  -11.35%  '
$ws.Range("D18").Value = '28.421.54'
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").Value = '''5.367'
$ws.Range("E19").Value = '  -5.68%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").Value = '''207.09'
$ws.Range("E21").Value = '  -6.16%  '
$ws.Range("D22").Value = '''11.34'
$ws.Range("E22").Value = '  -1.87%  '
$ws.Range("D23").Value = '''6.736'
$ws.Range("E23").Value = '  -1.72%  '
$ws.Range("D24").Value = '''1.002'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = '''151.79'
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("D26").Value = '''8.080'
$ws.Range("E26").Value = '  +1.52%  '
$ws.Range("D27").Value = '''0.1248'
$ws.Range("E27").Value = '  -3.19%  '
$ws.Range("E28").Value = '  -1.95%  '
$ws.Range("D29").Value = '''1.407'
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("D30").Value = '''0.06145'
$ws.Range("E30").Value = '  -4.27%  '
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("D32").Value = '''3.759'
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").Value = '''3.750'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").Value = '''1.669'
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("D35").Value = '''1.043'
$ws.Range("E35").Value = '  -4.37%  '
$ws.Range("D36").Value = '''0.6372'
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").Value = '''2.500'
$ws.Range("E37").Value = '  -1.40%  '
$ws.Range("D38").Value = '''2.678'
$ws.Range("E38").Value = '  -2.01%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.142.68'
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '''0.01668'
$ws.Range("E40").Value = '  -4.71%  '
$ws.Range("E41").Value = '  -4.55%  '
$ws.Range("D42").Value = '''0.8744'
$ws.Range("E42").Value = '  -2.09%  '
$ws.Range("D43").Value = '''1.007'
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").Value = '''99.80'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").Value = '1.932.49'
$ws.Range("D46").Value = '''59.61'
$ws.Range("E46").Value = '  -4.01%  '
$ws.Range("D47").Value = '''0.00000000110'
$ws.Range("E47").Value = '  -2.45%  '
$ws.Range("D48").Value = '''8.398'
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").Value = '''1.568'
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("D50").Value = '''0.05408'
$ws.Range("E50").Value = '  -3.18%  '
$ws.Range("D51").Value = '''0.4465'
$ws.Range("E51").Value = '  -1.92%  '
